# Update odds values in the "Jogos da Semana" worksheet to reflect the
# latest FlashScore snapshot for 2024-11-26.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Mumbai City vs Punjab) ---
$ws.Range("G2").Value  = 1.75
$ws.Range("I2").Value  = 4.2
$ws.Range("J2").Value  = 2.38
$ws.Range("L2").Value  = 4.33
$ws.Range("U2").Value  = 1.62
$ws.Range("V2").Value  = 2.2
$ws.Range("X2").Value  = 9.5
$ws.Range("AK2").Value = 41
$ws.Range("AL2").Value = 29
$ws.Range("AW2").Value = 6
$ws.Range("AX2").Value = 21

# --- Row 3 (River Plate vs CA Cerro) ---
$ws.Range("G3").Value  = 2.38
$ws.Range("I3").Value  = 3.25
$ws.Range("J3").Value  = 3.1
$ws.Range("N3").Value  = 8
$ws.Range("Q3").Value  = 2.2
$ws.Range("R3").Value  = 1.65
$ws.Range("W3").Value  = 7
$ws.Range("AO3").Value = 13
$ws.Range("AR3").Value = 67
$ws.Range("AW3").Value = 5
